# Updated RAD EL-Motor Fuel Tax.
# Refresh the "Date" column (column B) execution timestamps on Sheet1
# to reflect the latest Katalon test run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Mon Apr 08 18:06:07 EDT 2024"
$ws.Range("B3").Value = "Mon Apr 08 18:06:18 EDT 2024"
$ws.Range("B4").Value = "Mon Apr 08 18:06:30 EDT 2024"
$ws.Range("B5").Value = "Mon Apr 08 18:06:41 EDT 2024"
$ws.Range("B6").Value = "Mon Apr 08 18:06:53 EDT 2024"
$ws.Range("B7").Value = "Mon Apr 08 18:07:04 EDT 2024"
